# Updates the "cryptos" price-list sheet (GitHub Actions scrape refresh):
# new Price (col D) / Volume(1h) (col E) readings for most coins, and the
# dogwifhat / FirstDigitalUSD rows (50-51) swap order with updated figures.
#
# All cells in this sheet are plain text (prices like "64.487.81" or
# "1.00" must stay text, not become numbers), so every write goes through
# Set-TextValue: force NumberFormat to "@" (Text) before assigning the
# value so Excel doesn't auto-convert numeric-looking strings, then put
# the cell's style back to "Normal" so no stray number formatting is left
# behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-TextValue "D2" "64.487.81"
Set-TextValue "E2" "  -2.79%  "
Set-TextValue "D3" "3.178.88"
Set-TextValue "E3" "  -4.31%  "
Set-TextValue "E4" "  +0.02%  "
Set-TextValue "D5" "569.57"
Set-TextValue "E5" "  -3.17%  "
Set-TextValue "E6" "  -7.95%  "
Set-TextValue "D7" "0.607"
Set-TextValue "E7" "  -6.22%  "
Set-TextValue "D9" "3.179.65"
Set-TextValue "E9" "  -4.28%  "
Set-TextValue "E10" "  -4.04%  "
Set-TextValue "E11" "  +0.14%  "
Set-TextValue "E12" "  -3.72%  "
Set-TextValue "D13" "3.731.35"
Set-TextValue "E13" "  -4.42%  "
Set-TextValue "E14" "  -2.33%  "
Set-TextValue "D15" "64.547.23"
Set-TextValue "E15" "  -2.73%  "
Set-TextValue "D16" "25.39"
Set-TextValue "E16" "  -3.18%  "
Set-TextValue "E17" "  -2.67%  "
Set-TextValue "D18" "3.153.01"
Set-TextValue "E18" "  -4.94%  "
Set-TextValue "D19" "419.58"
Set-TextValue "E19" "  -1.64%  "
Set-TextValue "E20" "  -2.39%  "
Set-TextValue "D21" "5.36"
Set-TextValue "E21" "  -3.15%  "
Set-TextValue "D22" "7.11"
Set-TextValue "E22" "  -4.10%  "
Set-TextValue "D23" "0.999"
Set-TextValue "E23" "  -0.25%  "
Set-TextValue "E24" "  -0.39%  "
Set-TextValue "D25" "70.12"
Set-TextValue "D26" "0.205"
Set-TextValue "E26" "  +0.27%  "
Set-TextValue "D27" "0.488"
Set-TextValue "E27" "  -5.32%  "
Set-TextValue "E28" "  -7.11%  "
Set-TextValue "D29" "8.86"
Set-TextValue "E29" "  -1.08%  "
Set-TextValue "D30" "1.00"
Set-TextValue "E30" "  +0.09%  "
Set-TextValue "E31" "  -5.58%  "
Set-TextValue "D32" "21.71"
Set-TextValue "E32" "  -2.85%  "
Set-TextValue "E33" "  -0.10%  "
Set-TextValue "E34" "  -2.46%  "
Set-TextValue "D35" "6.33"
Set-TextValue "E35" "  -4.03%  "
Set-TextValue "E36" "  -4.33%  "
Set-TextValue "D37" "157.65"
Set-TextValue "E37" "  -1.31%  "
Set-TextValue "D39" "2.726.23"
Set-TextValue "E39" "  -5.69%  "
Set-TextValue "E40" "  -5.11%  "
Set-TextValue "D41" "24.30"
Set-TextValue "E41" "  -8.28%  "
Set-TextValue "D42" "4.18"
Set-TextValue "E42" "  -3.23%  "
Set-TextValue "D43" "39.18"
Set-TextValue "E43" "  -2.43%  "
Set-TextValue "E45" "  -6.26%  "
Set-TextValue "D46" "5.62"
Set-TextValue "E46" "  -5.18%  "
Set-TextValue "E47" "  -2.93%  "
Set-TextValue "D48" "294.29"
Set-TextValue "E48" "  -6.41%  "
Set-TextValue "E49" "  -6.88%  "
Set-TextValue "B50" "FirstDigitalUSD"
Set-TextValue "C50" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D50" "1.00"
Set-TextValue "E50" "  +0.02%  "
Set-TextValue "B51" "dogwifhat"
Set-TextValue "C51" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D51" "2.00"
Set-TextValue "E51" "  -13.27%  "
